$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Macroeconomics B" entry from F5 to E6.
$srcCell = $ws.Range("F5")
$dstCell = $ws.Range("E6")

# Copy value + formatting (fill, font, alignment) from F5 into E6.
$srcCell.Copy($dstCell)

# Remove the now-vacated source cell.
$srcCell.Clear()

# Update the active selection to F6, matching the saved workbook view.
$ws.Range("F6").Select()
